$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the source data (RM 232 and SC 92).
# Delete the higher-numbered row first so the lower row index stays valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# Apply the updated / cleared values for the remaining rows (post row-shift numbering).
$ws.Range("D2").Value = -13.5
$ws.Range("F3").Value = ""
$ws.Range("F4").Value = 17.97
$ws.Range("D6").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("D12").Value = -14.1
$ws.Range("D14").Value = ""
$ws.Range("F15").Value = 16.2
$ws.Range("F18").Value = 18.35
$ws.Range("F19").Value = ""
$ws.Range("D20").Value = -14
$ws.Range("D21").Value = -14.3
$ws.Range("F22").Value = ""
$ws.Range("D23").Value = ""
$ws.Range("F23").Value = 16.48
$ws.Range("D24").Value = ""
$ws.Range("F25").Value = 16.6
$ws.Range("C26").Value = 10.8
$ws.Range("C27").Value = ""
$ws.Range("F27").Value = ""
$ws.Range("C28").Value = ""
$ws.Range("C29").Value = 11.2
$ws.Range("C30").Value = 11.4
$ws.Range("C31").Value = ""
$ws.Range("D31").Value = -13.7
$ws.Range("C32").Value = ""
$ws.Range("D33").Value = -14.1
